$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "37.504.95"
$ws.Range("E2").Value = "  +2.58%  "
$ws.Range("D3").Value = "2.077.48"
$ws.Range("E3").Value = "  +3.59%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = "'235.15"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.01%  "
$ws.Range("D6").Value = "'0.623"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +3.58%  "
$ws.Range("D7").Value = "'58.17"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +5.65%  "
$ws.Range("E8").Value = "  +0.00%  "
$ws.Range("D9").Value = "'0.386"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +3.79%  "
$ws.Range("D10").Value = "'59.08"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.27%  "
$ws.Range("E11").Value = "  +1.88%  "
$ws.Range("D12").Value = "'0.101"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +3.35%  "
$ws.Range("D13").Value = "2.383.00"
$ws.Range("E13").Value = "  +3.65%  "
$ws.Range("E14").Value = "  +2.24%  "
$ws.Range("D15").Value = "'21.14"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +4.09%  "
$ws.Range("D16").Value = "'0.780"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +2.83%  "
$ws.Range("E17").Value = "  +1.59%  "
$ws.Range("D18").Value = "2.075.57"
$ws.Range("E18").Value = "  +3.48%  "
$ws.Range("D19").Value = "37.640.04"
$ws.Range("E19").Value = "  +3.07%  "
$ws.Range("E20").Value = "  +17.86%  "
$ws.Range("D21").Value = "'69.90"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +2.99%  "
$ws.Range("D22").Value = "0.0₃0817"
$ws.Range("E22").Value = "  +1.48%  "
$ws.Range("D23").Value = "'226.47"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +2.05%  "
$ws.Range("D24").Value = "'1.00"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.06%  "
$ws.Range("D25").Value = "'2.49"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +3.34%  "
$ws.Range("D26").Value = "'2.39"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.47%  "
$ws.Range("D27").Value = "'167.34"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +2.86%  "
$ws.Range("E28").Value = "  +10.04%  "
$ws.Range("D29").Value = "'9.01"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +3.83%  "
$ws.Range("D30").Value = "'19.30"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +2.42%  "
$ws.Range("E31").Value = "  -1.41%  "
$ws.Range("E32").Value = "  +1.10%  "
$ws.Range("D33").Value = "'4.52"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +2.79%  "
$ws.Range("B34").Value = "Hedera"
$ws.Range("C34").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D34").Value = "'0.0625"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +3.09%  "
$ws.Range("B35").Value = "LidoDAOToken"
$ws.Range("C35").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D35").Value = "'2.58"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +6.03%  "
$ws.Range("D36").Value = "'4.60"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +7.63%  "
$ws.Range("E37").Value = "  -0.06%  "
$ws.Range("D38").Value = "'3.34"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.38%  "
$ws.Range("E39").Value = "  +2.68%  "
$ws.Range("D40").Value = "'1.76"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.05%  "
$ws.Range("D41").Value = "'4.63"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +21.04%  "
$ws.Range("E42").Value = "  -1.06%  "
$ws.Range("D43").Value = "'0.0962"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +4.01%  "
$ws.Range("D44").Value = "1.476.64"
$ws.Range("E44").Value = "  +1.61%  "
$ws.Range("E45").Value = "  +6.95%  "
$ws.Range("D46").Value = "'95.88"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +5.80%  "
$ws.Range("E47").Value = "  +4.60%  "
$ws.Range("D48").Value = "'15.82"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +3.13%  "
$ws.Range("E49").Value = "  +3.56%  "
$ws.Range("D50").Value = "'7.28"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +5.94%  "
$ws.Range("E51").Value = "  +1.83%  "
